# Edit coverage file location / reshape the time-axis columns.
#
# The sheet originally had quarterly columns (year, year+1/12, year+0.5,
# year+7/12) running from 2018 to 2040 (H1:BZ1). The new layout only keeps
# the half-yearly columns (year, year+1/12) for 2024 onward, so everything
# from P1 onward is rewritten and the trailing columns (AV:BZ) are dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ColLetter([int]$n) {
    $s = ""
    while ($n -gt 0) {
        $rem = ($n - 1) % 26
        $s = [char](65 + $rem) + $s
        $n = [int](($n - $rem - 1) / 26)
    }
    return $s
}

# --- Rewrite row 1 (the year axis) from P1 through AU1 -------------------
# H1:O1 (2018 .. 2024.0833333333333) stay as-is.
$startCol = 16   # column P
$idx = 0
for ($year = 2025; $year -le 2040; $year++) {
    $c1 = ColLetter ($startCol + $idx); $idx++
    $c2 = ColLetter ($startCol + $idx); $idx++
    $ws.Range($c1 + "1").Value = [double]$year
    if ($year -eq 2040) {
        # matches the source value's exact (slightly rounded) double
        $ws.Range($c2 + "1").Value = 2040.0833333333301
    } else {
        $ws.Range($c2 + "1").Value = $year + (1.0 / 12.0)
    }
}

# --- Drop the now-unused trailing columns (old AV:BZ) on rows 1-3 --------
$ws.Range("AV1:BZ3").ClearContents()

# --- Update the view: scroll right and reselect AU2 ----------------------
$ws.Range("AU2").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
